$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column before A; this shifts the existing Task/Start Date/
#    Frequency/Day of Month/Heure columns from A:E to B:F (formatting, column
#    widths and styles all travel with the shift automatically).
# ---------------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# ---------------------------------------------------------------------------
# 2) New "GROUP" column (A). Header is bold (copy header look from B1 then
#    tweak), data cells are plain default-styled text.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "GROUP"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2:A9").Value = "CRM"

# ---------------------------------------------------------------------------
# 3) Small text / value fixes on the existing (shifted) columns.
# ---------------------------------------------------------------------------
# typo fix
$ws.Range("B3").Value = "Welcome Journey"
# date correction for "Product reco"
$ws.Range("C5").Value = 45720
# header rename Heure -> Hour
$ws.Range("F1").Value = "Hour"

# ---------------------------------------------------------------------------
# 4) F6:F9 used to carry the old stand-alone "time, no wrap" style; it should
#    now look like the rest of the Hour column (wrap + vertical-center, same
#    as F2). Copy the format straight from F2 so the stylesheet entry is
#    reused instead of minted anew.
# ---------------------------------------------------------------------------
$ws.Range("F2").Copy()
$ws.Range("F6:F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 7:9 no longer need the extra wrapped row height now that column B is
# wider - let Excel recompute their natural height.
$ws.Rows("7:9").AutoFit()

# ---------------------------------------------------------------------------
# 5) Two new tasks/rows for the INGESTION group, formatted like row 9.
# ---------------------------------------------------------------------------
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = "INGESTION"
$ws.Range("B10").Value = "Segmentation"
$ws.Range("C10").Value = 45717
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = 0.875

$ws.Range("A11").Value = "INGESTION"
$ws.Range("B11").Value = "Recommendation"
$ws.Range("C11").Value = 45717
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 0.77083333333333337

# ---------------------------------------------------------------------------
# 6) Column widths: column A reverts to the sheet default (no custom width);
#    column B (now holding the Task names) needs to be wide enough for
#    "Product repurchase" etc.
# ---------------------------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 21

# ---------------------------------------------------------------------------
# 7) Selection, matching the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("C5").Select()
